$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coding Progress")

# Fix the missing "partial DONE" marker for the "Master Group" block (row 3)
$ws.Range("C3").Value = "partial DONE"

# Add a new "Master Category" block (rows 14-16), mirroring the pattern
# used by the other Master* blocks (e.g. Master User rows 5-7)
$ws.Range("B14").Value = "Master Category"
$ws.Range("C14").Value = "DONE"
$ws.Range("D14").Value = "Data Entry Coding"

$ws.Range("C15").Value = "partial DONE"
$ws.Range("D15").Value = "Input validation"

$ws.Range("C16").Value = "DONE"
$ws.Range("D16").Value = "DataGrid for browsing"

# Match the saved selection state from the diff
$ws.Range("C3").Select()
